$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume table refresh (GitHub Actions data pull).
# Column D holds prices as text (some rows use "." as a thousands
# separator, e.g. "65.764.23"), so force text format before writing
# any value that Excel would otherwise auto-parse as a number.

$ws.Range("D2").Value = "65.764.23"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "2.656.06"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.42"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.99"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.628"
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("B10").Value = "Cardano"
$ws.Range("C10").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.395"
$ws.Range("E10").Value = "  -1.36%  "
$ws.Range("B11").Value = "Toncoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.80"
$ws.Range("E11").Value = "  -0.95%  "
$ws.Range("E12").Value = "  +1.32%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.51"
$ws.Range("E13").Value = "  -2.56%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000196"
$ws.Range("E14").Value = "  +0.59%  "
$ws.Range("D15").Value = "3.130.06"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").Value = "65.578.99"
$ws.Range("D17").Value = "2.614.43"
$ws.Range("E17").Value = "  -1.28%  "
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.45"
$ws.Range("E20").Value = "  -1.78%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.12"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("E25").Value = "  +6.93%  "
$ws.Range("E26").Value = "  -0.90%  "
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.164"
$ws.Range("E28").Value = "  -1.91%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "544.56"
$ws.Range("E29").Value = "  +2.09%  "
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.90"
$ws.Range("E31").Value = "  -1.73%  "
$ws.Range("E32").Value = "  +0.54%  "
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.46"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.42"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.34"
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.94"
$ws.Range("E39").Value = "  +0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.82"
$ws.Range("E40").Value = "  -2.79%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "161.50"
$ws.Range("E42").Value = "  -1.66%  "
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0603"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.55"
$ws.Range("E46").Value = "  -1.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.637"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0998"
$ws.Range("E49").Value = "  -1.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.72"
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("E51").Value = "  +7.04%  "
